# Populate columns B-E for rows 209-301 ("Add 100 more items")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple: RowNumber, B(category text), C(value - numeric unless CIsText), CIsText, D(text), E(text)
$rowsData = @(
    @(209, 'Light Armor', '2', $false, 'Combat', 'Uncommon'),
    @(210, 'Light Armor', '3', $false, 'Utility', 'Rare'),
    @(211, 'Light Armor', '1', $false, 'Utility', 'Very Rare'),
    @(212, 'Light Armor', '1', $false, 'Combat', 'Uncommon'),
    @(213, 'Light Armor', '2', $false, 'Utility', 'Rare'),
    @(214, 'Light Armor', '2', $false, 'Combat', 'Very Rare'),
    @(215, 'Light Armor', '1', $false, 'Utility', 'Rare'),
    @(216, 'Light Armor', '1', $false, 'Utility', 'Uncommon'),
    @(217, 'Light Armor', '1', $false, 'Utility', 'Uncommon'),
    @(218, 'Light Armor', '1', $false, 'Utility', 'Rare'),
    @(219, 'Light Armor', '1', $false, 'Combat', 'Uncommon'),
    @(220, 'Light Armor', '2', $false, 'Combat', 'Very Rare'),
    @(221, 'Light Armor', '2', $false, 'Utility', 'Very Rare'),
    @(222, 'Heavy Armor', '2', $false, 'Utility', 'Rare'),
    @(223, 'Heavy Armor', '2', $false, 'Combat', 'Rare'),
    @(224, 'Heavy Armor', '2', $false, 'Combat', 'Rare'),
    @(225, 'Heavy Armor', '1', $false, 'Combat', 'Rare'),
    @(226, 'Heavy Armor', '1', $false, 'Combat', 'Rare'),
    @(227, 'Heavy Armor', '2', $false, 'Combat', 'Legendary'),
    @(228, 'Heavy Armor', '1', $false, 'Combat', 'Very Rare'),
    @(229, 'Heavy Armor', '2', $false, 'Combat', 'Very Rare'),
    @(230, 'Heavy Armor', '1', $false, 'Combat', 'Very Rare'),
    @(231, 'Heavy Armor', '1', $false, 'Combat', 'Rare'),
    @(232, 'Heavy Armor', '3', $false, 'Combat', 'Rare'),
    @(233, 'Heavy Armor', '1', $false, 'Combat', 'Rare'),
    @(234, 'Heavy Armor', '3', $false, 'Combat', 'Very Rare'),
    @(235, 'Heavy Armor', '1', $false, 'Utility', 'Rare'),
    @(236, 'Heavy Armor', '1', $false, 'Combat', 'Legendary'),
    @(237, 'Heavy Armor', '2', $false, 'Combat', 'Very Rare'),
    @(238, 'Heavy Armor', '2', $false, 'Combat', 'Rare'),
    @(239, 'Heavy Armor', '2', $false, 'Utility', 'Rare'),
    @(240, 'Heavy Armor', '1', $false, 'Utility', 'Rare'),
    @(241, 'Heavy Armor', '2', $false, 'Combat', 'Uncommon'),
    @(242, 'Helmet', '1', $false, 'Utility', 'Rare'),
    @(243, 'Helmet', '1', $false, 'Combat', 'Uncommon'),
    @(244, 'Helmet', '1', $false, 'Combat', 'Rare'),
    @(245, 'Helmet', '1', $false, 'Utility', 'Rare'),
    @(246, 'Helmet', '4', $false, 'Cursed', 'Rare'),
    @(247, 'Helmet', '3', $false, 'Combat', 'Rare'),
    @(248, 'Helmet', '1', $false, 'Utility', 'Uncommon'),
    @(249, 'Helmet', '1', $false, 'Utility', 'Uncommon'),
    @(250, 'Helmet', '2', $false, 'Cursed', 'Very Rare'),
    @(251, 'Helmet', '1', $false, 'Utility', 'Rare'),
    @(252, 'Cuirass/Breastplate', '1', $false, 'Utility', 'Rare'),
    @(253, 'Cuirass/Breastplate', '1', $false, 'Utility', 'Legendary'),
    @(254, 'Cuirass/Breastplate', '2', $false, 'Cursed', 'Rare'),
    @(255, 'Cuirass/Breastplate', '2', $false, 'Combat', 'Legendary'),
    @(256, 'Cuirass/Breastplate', '3', $false, 'Combat', 'Rare'),
    @(257, 'Cuirass/Breastplate', '1', $false, 'Utility', 'Uncommon'),
    @(258, 'Cuirass/Breastplate', '2', $false, 'Combat', 'Rare'),
    @(259, 'Cuirass/Breastplate', '2', $false, 'Utility', 'Rare'),
    @(260, 'Cuirass/Breastplate', '1', $false, 'Combat', 'Legendary'),
    @(261, 'Cuirass/Breastplate', '1', $false, 'Utility', 'Very Rare'),
    @(262, 'Bracer/Gauntlet', '3', $false, 'Utility', 'Very Rare'),
    @(263, 'Bracer/Gauntlet', '1', $false, 'Utility', 'Uncommon'),
    @(264, 'Bracer/Gauntlet', '1', $false, 'Combat', 'Uncommon'),
    @(265, 'Bracer/Gauntlet', '2', $false, 'Combat', 'Rare'),
    @(266, 'Bracer/Gauntlet', '1', $false, 'Utility', 'Rare'),
    @(267, 'Bracer/Gauntlet', '1', $false, 'Combat', 'Legendary'),
    @(268, 'Bracer/Gauntlet', '2', $false, 'Combat', 'Uncommon'),
    @(269, 'Bracer/Gauntlet', '1', $false, 'Combat', 'Uncommon'),
    @(270, 'Bracer/Gauntlet', '3', $false, 'Combat', 'Uncommon'),
    @(271, 'Bracer/Gauntlet', '1', $false, 'Combat', 'Uncommon'),
    @(272, 'Bracer/Gauntlet', '1', $false, 'Combat', 'Rare'),
    @(273, 'Bracer/Gauntlet', '2', $false, 'Utility', 'Rare'),
    @(274, 'Bracer/Gauntlet', '3', $false, 'Combat', 'Rare'),
    @(275, 'Bracer/Gauntlet', '1', $false, 'Utility', 'Rare'),
    @(276, 'Bracer/Gauntlet', '1', $false, 'Utility', 'Very Rare'),
    @(277, 'Bracer/Gauntlet', '1', $false, 'Combat', 'Rare'),
    @(278, 'NA', 'NA', $true, 'NA', 'NA'),
    @(279, 'Bracer/Gauntlet', '1', $false, 'Combat', 'Very Rare'),
    @(280, 'Bracer/Gauntlet', '1', $false, 'Combat', 'Very Rare'),
    @(281, 'Bracer/Gauntlet', '2', $false, 'Combat', 'Very Rare'),
    @(282, 'Greaves', '1', $false, 'Utility', 'Rare'),
    @(283, 'Greaves', '1', $false, 'Utility', 'Uncommon'),
    @(284, 'Greaves', '2', $false, 'Cursed', 'Uncommon'),
    @(285, 'Greaves', '2', $false, 'Combat', 'Uncommon'),
    @(286, 'Greaves', '1', $false, 'Combat', 'Uncommon'),
    @(287, 'Greaves', '1', $false, 'Utility', 'Uncommon'),
    @(288, 'Greaves', '1', $false, 'Utility', 'Uncommon'),
    @(289, 'Greaves', '1', $false, 'Utility', 'Uncommon'),
    @(290, 'Greaves', '2', $false, 'Utility', 'Rare'),
    @(291, 'Greaves', '1', $false, 'Utility', 'Common'),
    @(292, 'Shields', '1', $false, 'Combat', 'Very Rare'),
    @(293, 'Shields', '2', $false, 'Combat', 'Legendary'),
    @(294, 'Shields', '3', $false, 'Combat', 'Very Rare'),
    @(295, 'Shields', '2', $false, 'Combat', 'Uncommon'),
    @(296, 'Shields', '1', $false, 'Utility', 'Uncommon'),
    @(297, 'Shields', '1', $false, 'Combat', 'Rare'),
    @(298, 'Shields', '1', $false, 'Utility', 'Uncommon'),
    @(299, 'Shields', '2', $false, 'Combat', 'Rare'),
    @(300, 'Shields', '1', $false, 'Combat', 'Uncommon'),
    @(301, 'Shields', '2', $false, 'Combat', 'Uncommon')
)

foreach ($item in $rowsData) {
    $rowNum   = $item[0]
    $bText    = $item[1]
    $cValue   = $item[2]
    $cIsText  = $item[3]
    $dText    = $item[4]
    $eText    = $item[5]

    $ws.Cells.Item($rowNum, 2).Value = $bText
    if ($cIsText) {
        $ws.Cells.Item($rowNum, 3).Value = $cValue
    } else {
        $ws.Cells.Item($rowNum, 3).Value = [double]$cValue
    }
    $ws.Cells.Item($rowNum, 4).Value = $dText
    $ws.Cells.Item($rowNum, 5).Value = $eText
}

# Reflect the cursor/selection position left behind after entering the new rows
[void]$ws.Range("H290").Select()
